# Updated cryptos list on Fri Mar 29 18:43:27 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay TEXT even when it looks numeric
# (Price column uses dot-grouped / fixed-decimal strings such as "610.00"
# or "69.452.49" that must not be coerced into Excel numbers).
function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

$ws.Range("D2").Value = "69.452.49"
$ws.Range("E2").Value = "  -1.94%  "
$ws.Range("D3").Value = "3.491.23"
$ws.Range("E3").Value = "  -1.88%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.16%  "
Set-TextValue $ws.Range("D5") "610.00"
$ws.Range("E5").Value = "  +4.74%  "
Set-TextValue $ws.Range("D6") "186.03"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("E8").Value = "  -0.16%  "
Set-TextValue $ws.Range("D9") "0.215"
$ws.Range("E9").Value = "  -2.16%  "
Set-TextValue $ws.Range("D10") "0.651"
$ws.Range("E10").Value = "  -0.04%  "
Set-TextValue $ws.Range("D11") "53.12"
$ws.Range("E11").Value = "  -2.44%  "
$ws.Range("E12").Value = "  -2.42%  "
Set-TextValue $ws.Range("D13") "9.53"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "4.035.00"
$ws.Range("E14").Value = "  -2.19%  "
Set-TextValue $ws.Range("D15") "603.90"
$ws.Range("E15").Value = "  +6.03%  "
$ws.Range("D16").Value = "69.485.98"
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("E17").Value = "  +1.27%  "
Set-TextValue $ws.Range("D18") "18.88"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "3.490.00"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("E21").Value = "  -1.29%  "
Set-TextValue $ws.Range("D22") "17.14"
$ws.Range("E22").Value = "  -2.77%  "
Set-TextValue $ws.Range("D23") "105.89"
$ws.Range("E23").Value = "  +11.78%  "
$ws.Range("E24").Value = "  +1.90%  "
Set-TextValue $ws.Range("D25") "5.04"
$ws.Range("E25").Value = "  +2.86%  "
$ws.Range("E26").Value = "  +2.79%  "
Set-TextValue $ws.Range("D27") "10.94"
$ws.Range("E27").Value = "  -2.57%  "
Set-TextValue $ws.Range("D28") "9.74"
$ws.Range("E28").Value = "  +6.42%  "
Set-TextValue $ws.Range("D29") "33.64"
$ws.Range("E29").Value = "  +3.65%  "
Set-TextValue $ws.Range("D30") "6.98"
$ws.Range("E30").Value = "  -3.10%  "
Set-TextValue $ws.Range("D31") "12.44"
$ws.Range("E31").Value = "  +1.35%  "
Set-TextValue $ws.Range("D32") "4.06"
$ws.Range("E32").Value = "  +20.07%  "
$ws.Range("E33").Value = "  -0.68%  "
Set-TextValue $ws.Range("D34") "63.25"
$ws.Range("E34").Value = "  +0.37%  "
Set-TextValue $ws.Range("D35") "3.19"
$ws.Range("E35").Value = "  -6.61%  "
$ws.Range("E36").Value = "  -0.19%  "
Set-TextValue $ws.Range("D37") "525.04"
$ws.Range("E37").Value = "  -3.83%  "
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("D39").Value = "3.611.21"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("E40").Value = "  +5.64%  "
Set-TextValue $ws.Range("D41") "36.76"
$ws.Range("E41").Value = "  -2.96%  "
$ws.Range("D42").Value = "0.0₃0779"
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("E44").Value = "  -0.99%  "
Set-TextValue $ws.Range("D45") "2.94"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("E47").Value = "  -4.58%  "
$ws.Range("E48").Value = "  -5.52%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("E50").Value = "  -9.53%  "
Set-TextValue $ws.Range("D51") "0.000243"
$ws.Range("E51").Value = "  -8.10%  "
